# Correções para Dia 2
# Set cells D4, F4, D5 and F5 on the "horario" sheet to the same value as
# B4/B5 ("Esporte"), and update the active selection to F4:F5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("horario")

$value = $ws.Range("B4").Value()

$ws.Range("D4").Value = $value
$ws.Range("F4").Value = $value
$ws.Range("D5").Value = $value
$ws.Range("F5").Value = $value

$ws.Activate()
$ws.Range("F4:F5").Select()

$wsCarga = $wb.Worksheets.Item("carga")
$wsCarga.Activate()
$wsCarga.Range("F4:F5,A11").Select()

$wsDias = $wb.Worksheets.Item("dias-atividade")
$wsDias.Activate()
$wsDias.Range("F4:F5,H9").Select()

$ws.Activate()
